$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value2 = 114573.89
$ws.Range("I40").Value2 = 751530
$ws.Range("J40").Value2 = 3798.913
$ws.Range("K40").Value2 = 751530
$ws.Range("L40").Value2 = 3798.913
$ws.Range("M40").Value2 = -751355
$ws.Range("N40").Value2 = -4148.913
$ws.Range("H62").Value2 = 8339.786
$ws.Range("I62").Value2 = 6998.375
$ws.Range("K62").Value2 = 6998.375
$ws.Range("M62").Value2 = -6374.375
$ws.Range("H65").Value2 = 8339.786
$ws.Range("I65").Value2 = 6998.375
$ws.Range("K65").Value2 = 34991.875
$ws.Range("M65").Value2 = -31871.875
$ws.Range("H132").Value2 = 3725.9546
$ws.Range("I132").Value2 = 3787.9736
$ws.Range("K132").Value2 = 11363.9208
$ws.Range("M132").Value2 = -8833.9208
$ws.Range("H135").Value2 = 2041.0834
$ws.Range("J135").Value2 = 3099
$ws.Range("L135").Value2 = 27891
$ws.Range("N135").Value2 = -32961
$ws.Range("H141").Value2 = 4795.0713
$ws.Range("J141").Value2 = 8500
$ws.Range("L141").Value2 = 25500
$ws.Range("N141").Value2 = -35860

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 3614
$ws.Range("I74").Value2 = 2997
$ws.Range("K74").Value2 = 2997
$ws.Range("M74").Value2 = -2123
$ws.Range("H77").Value2 = 3614
$ws.Range("I77").Value2 = 2997
$ws.Range("K77").Value2 = 14985
$ws.Range("M77").Value2 = -10617
$ws.Range("H80").Value2 = 0
$ws.Range("J80").Value2 = 0
$ws.Range("L80").Value2 = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value2 = 0
$ws.Range("J83").Value2 = 0
$ws.Range("L83").Value2 = 0
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value2 = 1043.0714
$ws.Range("I97").Value2 = 718.2083
$ws.Range("K97").Value2 = 718.2083
$ws.Range("M97").Value2 = -222.2083
$ws.Range("H102").Value2 = 3123.8
$ws.Range("I102").Value2 = 2939.6667
$ws.Range("J102").Value2 = 3400
$ws.Range("K102").Value2 = 2939.6667
$ws.Range("L102").Value2 = 3400
$ws.Range("M102").Value2 = -1317.6667
$ws.Range("N102").Value2 = -6644
$ws.Range("H110").Value2 = 2242.2856
$ws.Range("I110").Value2 = 1565.3334
$ws.Range("K110").Value2 = 1565.3334
$ws.Range("M110").Value2 = 479.6666
$ws.Range("H138").Value2 = 97249
$ws.Range("J138").Value2 = 97249
$ws.Range("L138").Value2 = 97249
$ws.Range("N138").Value2 = -107529

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value2 = 6674.25
$ws.Range("I5").Value2 = 6674.25
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = 6674.25
$ws.Range("L5").Value2 = 0
$ws.Range("M5").Value2 = -6561.25
$ws.Range("N5").ClearContents()
$ws.Range("H94").Value2 = 1458
$ws.Range("I94").Value2 = 687.9091
$ws.Range("J94").Value2 = 4281.6665
$ws.Range("K94").Value2 = 687.9091
$ws.Range("L94").Value2 = 4281.6665
$ws.Range("M94").Value2 = -236.9091
$ws.Range("N94").Value2 = -5183.6665
$ws.Range("H99").Value2 = 2600.889
$ws.Range("I99").Value2 = 2372.5715
$ws.Range("K99").Value2 = 2372.5715
$ws.Range("M99").Value2 = -874.5715
$ws.Range("H105").Value2 = 3502.1765
$ws.Range("I105").Value2 = 2198
$ws.Range("K105").Value2 = 2198
$ws.Range("M105").Value2 = -451
$ws.Range("H134").Value2 = 14994063
$ws.Range("I134").Value2 = 3761680.2
$ws.Range("K134").Value2 = 11285040.6
$ws.Range("M134").Value2 = -11282505.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value2 = 3795.8
$ws.Range("I2").Value2 = 3795.8
$ws.Range("K2").Value2 = 3795.8
$ws.Range("M2").Value2 = -3682.8
$ws.Range("H16").Value2 = 1805.9231
$ws.Range("I16").Value2 = 1330.7778
$ws.Range("K16").Value2 = 1330.7778
$ws.Range("M16").Value2 = -1043.7778
$ws.Range("H31").Value2 = 4001.1333
$ws.Range("I31").Value2 = 2210.5334
$ws.Range("J31").Value2 = 5791.7334
$ws.Range("K31").Value2 = 2210.5334
$ws.Range("L31").Value2 = 5791.7334
$ws.Range("M31").Value2 = -1915.5334
$ws.Range("N31").Value2 = -6381.7334
$ws.Range("H34").Value2 = 4001.1333
$ws.Range("I34").Value2 = 2210.5334
$ws.Range("J34").Value2 = 5791.7334
$ws.Range("K34").Value2 = 2210.5334
$ws.Range("L34").Value2 = 5791.7334
$ws.Range("M34").Value2 = -2008.5334
$ws.Range("N34").Value2 = -6195.7334
$ws.Range("H58").Value2 = 2176.879
$ws.Range("I58").Value2 = 1839.2069
$ws.Range("K58").Value2 = 1839.2069
$ws.Range("M58").Value2 = -1636.2069
$ws.Range("H99").Value2 = 1458.6666
$ws.Range("I99").Value2 = 1458.6666
$ws.Range("J99").Value2 = 0
$ws.Range("K99").Value2 = 1458.6666
$ws.Range("L99").Value2 = 0
$ws.Range("M99").Value2 = 39.33339999999998
$ws.Range("N99").ClearContents()
$ws.Range("H113").Value2 = 1805.9231
$ws.Range("I113").Value2 = 1330.7778
$ws.Range("K113").Value2 = 1330.7778
$ws.Range("M113").Value2 = 839.2221999999999
$ws.Range("H126").Value2 = 1458.6666
$ws.Range("I126").Value2 = 1458.6666
$ws.Range("J126").Value2 = 0
$ws.Range("K126").Value2 = 4375.9998
$ws.Range("L126").Value2 = 0
$ws.Range("M126").Value2 = -1905.9998
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value2 = 3322.4443
$ws.Range("I132").Value2 = 2843.1428
$ws.Range("K132").Value2 = 8529.428400000001
$ws.Range("M132").Value2 = -5999.428400000001
$ws.Range("H134").Value2 = 2574.6897
$ws.Range("I134").Value2 = 1348.6364
$ws.Range("K134").Value2 = 4045.9092
$ws.Range("M134").Value2 = -1510.9092
$ws.Range("H136").Value2 = 2176.879
$ws.Range("I136").Value2 = 1839.2069
$ws.Range("K136").Value2 = 5517.620699999999
$ws.Range("M136").Value2 = -2967.620699999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value2 = 844
$ws.Range("I70").Value2 = 844
$ws.Range("K70").Value2 = 2532
$ws.Range("M70").Value2 = -2217
$ws.Range("H73").Value2 = 844
$ws.Range("I73").Value2 = 844
$ws.Range("K73").Value2 = 2532
$ws.Range("M73").Value2 = -1440
$ws.Range("H129").Value2 = 1541.1818
$ws.Range("J129").Value2 = 1503.125
$ws.Range("L129").Value2 = 4509.375
$ws.Range("N129").Value2 = -14509.375
$ws.Range("H131").Value2 = 1688.5454
$ws.Range("I131").Value2 = 1197.3572
$ws.Range("J131").Value2 = 1856.2683
$ws.Range("K131").Value2 = 3592.0716
$ws.Range("L131").Value2 = 5568.8049
$ws.Range("M131").Value2 = 1447.9284
$ws.Range("N131").Value2 = -15648.8049

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 2509.889
$ws.Range("I80").Value2 = 1995
$ws.Range("K80").Value2 = 1995
$ws.Range("M80").Value2 = -997
$ws.Range("H83").Value2 = 2509.889
$ws.Range("I83").Value2 = 1995
$ws.Range("K83").Value2 = 9975
$ws.Range("M83").Value2 = -4983
$ws.Range("H97").Value2 = 905.125
$ws.Range("I97").Value2 = 698.6
$ws.Range("J97").Value2 = 1249.3334
$ws.Range("K97").Value2 = 698.6
$ws.Range("L97").Value2 = 1249.3334
$ws.Range("M97").Value2 = -202.6
$ws.Range("N97").Value2 = -2241.3334
$ws.Range("H122").Value2 = 1781.1904
$ws.Range("I122").Value2 = 1900.2307
$ws.Range("J122").Value2 = 1587.75
$ws.Range("K122").Value2 = 5700.6921
$ws.Range("L122").Value2 = 4763.25
$ws.Range("M122").Value2 = -3250.6921
$ws.Range("N122").Value2 = -9663.25
$ws.Range("H126").Value2 = 5514.5654
$ws.Range("I126").Value2 = 4657.6924
$ws.Range("J126").Value2 = 6628.5
$ws.Range("K126").Value2 = 13973.0772
$ws.Range("L126").Value2 = 19885.5
$ws.Range("M126").Value2 = -11503.0772
$ws.Range("N126").Value2 = -24825.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value2 = 20000
$ws.Range("I42").Value2 = 20000
$ws.Range("J42").Value2 = 0
$ws.Range("K42").Value2 = 20000
$ws.Range("L42").Value2 = 0
$ws.Range("M42").Value2 = -19437
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value2 = 20000
$ws.Range("I49").Value2 = 20000
$ws.Range("J49").Value2 = 0
$ws.Range("K49").Value2 = 20000
$ws.Range("L49").Value2 = 0
$ws.Range("M49").Value2 = -19853
$ws.Range("N49").ClearContents()
$ws.Range("H132").Value2 = 4938.5
$ws.Range("I132").Value2 = 4301.6
$ws.Range("K132").Value2 = 12904.8
$ws.Range("M132").Value2 = -10374.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value2 = 8578.444
$ws.Range("I126").Value2 = 8578.444
$ws.Range("K126").Value2 = 25735.332
$ws.Range("M126").Value2 = -23265.332
$ws.Range("H132").Value2 = 3489.913
$ws.Range("I132").Value2 = 2815.3157
$ws.Range("J132").Value2 = 6694.25
$ws.Range("K132").Value2 = 8445.947100000001
$ws.Range("L132").Value2 = 20082.75
$ws.Range("M132").Value2 = -5915.947100000001
$ws.Range("N132").Value2 = -25142.75
$ws.Range("H135").Value2 = 0
$ws.Range("J135").Value2 = 0
$ws.Range("L135").Value2 = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value2 = 17866
$ws.Range("I136").Value2 = 1347.1555
$ws.Range("J136").Value2 = 56989.58
$ws.Range("K136").Value2 = 4041.4665
$ws.Range("L136").Value2 = 170968.74
$ws.Range("M136").Value2 = -1491.4665
$ws.Range("N136").Value2 = -176068.74
